# Generate Report for Handoff
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# New column width (stored OOXML width 17.2159881591797 chars) expressed in the
# VBA "ColumnWidth" units the COM layer expects (stored = ColumnWidth + 5/6,
# quantized to 1/6 char) -> nearest achievable ColumnWidth.
$newColWidth = 16.38265482584637

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-11-09 01:10:11"

$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-11-09 01:09:58"

$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-11-09 01:10:11"

$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
